# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column G holds "K" values; update rows 2-10 with the newly computed values.
$kValues = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 0
    6  = 2
    7  = 3
    8  = 0
    9  = 4
    10 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
